# Refresh the cryptocurrency price/volume figures (scheduled data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D ("Price") cells are stored as plain text in the source data
# (e.g. "63.734.00", "4.40", "0.0560") rather than numbers, so the cell's
# number format is forced to Text before assignment. This keeps exact
# string formatting (trailing zeros, thousands "." separators, etc.)
# instead of letting Excel coerce numeric-looking strings into actual
# numbers, which would silently rewrite e.g. "4.40" -> 4.4.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.734.00'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.658.35'
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.65'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.13'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.69'
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.83'
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.130.69'
$ws.Range("E14").Value = '  +2.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.502.42'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.647.45'
$ws.Range("E17").Value = '  +2.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.42'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '343.99'
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.40'
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.77'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.04'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.71'
$ws.Range("E24").Value = '  +7.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.59'
$ws.Range("E25").Value = '  +11.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '565.75'
$ws.Range("E26").Value = '  +22.07%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.63'
$ws.Range("E28").Value = '  +3.86%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.96'
$ws.Range("E30").Value = '  +1.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.01'
$ws.Range("E31").Value = '  +4.26%  '
$ws.Range("E32").Value = '  +13.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0823'
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '175.31'
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.92'
$ws.Range("E35").Value = '  +9.17%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.27'
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  +6.22%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '169.83'
$ws.Range("E41").Value = '  +7.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.48'
$ws.Range("E42").Value = '  +2.81%  '
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.13'
$ws.Range("E44").Value = '  +4.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0560'
$ws.Range("E45").Value = '  +3.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.631'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0961'
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.92'
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.73'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.34'
$ws.Range("E51").Value = '  -0.95%  '
